# MLTracker.xlsx - "tidying up the repo"
#
# The "No Metro" sheet's evaluation-metric columns (B/C) are being
# repurposed: the header changes from r2/mse to roc_auc_score/average_precision,
# the first model's score row is updated to the new metric values, and all
# the other now-stale r2/mse numbers sprinkled through the table are cleared
# out (they no longer apply to the new metric pair).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("No Metro")

# --- Header row (row 1): rename the metric columns ------------------------
$ws.Range("B1").Value = "roc_auc_score"
$ws.Range("C1").Value = "average_precision"

# --- Row 2: refresh the first model's metric values ------------------------
$ws.Range("B2").Value = 0.52
$ws.Range("C2").Value = 0.71

# --- Remaining rows: drop the now-obsolete r2/mse figures -------------------
$staleRows = @(3, 5, 6, 7, 9, 10, 11, 14, 15, 17, 18, 19, 22, 23, 26, 27, 29, 30, 31, 32, 33)
foreach ($r in $staleRows) {
    $ws.Range("B$r`:C$r").Clear()
}

# --- Restore the view state left by the last person who had it open --------
[void]$ws.Range("A10").Select()
